# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Ciruela - Angeleno, 2023-03-28 / serial 45013)
# right before the existing row 707, shifting all subsequent rows down by 2
# (old row 707 -> new row 709, ..., old row 787 -> new row 789).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 707-708; everything from old row 707 onward shifts down.
$ws.Rows("707:708").Insert()

# Populate the newly-inserted row 707.
$ws.Range("A707").Value = 6
$ws.Range("B707").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C707").Value = "Metropolitana"
$ws.Range("D707").Value = 45013
$ws.Range("E707").Value = 13
$ws.Range("F707").Value = "Fruta"
$ws.Range("G707").Value = 100103
$ws.Range("H707").Value = "Frutos de hueso (carozo)"
$ws.Range("I707").Value = 100103002
$ws.Range("J707").Value = "Ciruela"
$ws.Range("K707").Value = "Angeleno"
$ws.Range("L707").Value = "Primera"
$ws.Range("M707").Value = 24
$ws.Range("N707").Value = 160000
$ws.Range("O707").Value = 170000
$ws.Range("P707").Value = 165000
$ws.Range("Q707").Value = "$/bins (450 kilos)"
$ws.Range("R707").Value = "Región de O'Higgins"
$ws.Range("S707").Value = 367
$ws.Range("T707").Value = 450

# Populate the newly-inserted row 708.
$ws.Range("A708").Value = 6
$ws.Range("B708").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C708").Value = "Metropolitana"
$ws.Range("D708").Value = 45013
$ws.Range("E708").Value = 13
$ws.Range("F708").Value = "Fruta"
$ws.Range("G708").Value = 100103
$ws.Range("H708").Value = "Frutos de hueso (carozo)"
$ws.Range("I708").Value = 100103002
$ws.Range("J708").Value = "Ciruela"
$ws.Range("K708").Value = "Angeleno"
$ws.Range("L708").Value = "Segunda"
$ws.Range("M708").Value = 19
$ws.Range("N708").Value = 130000
$ws.Range("O708").Value = 130000
$ws.Range("P708").Value = 130000
$ws.Range("Q708").Value = "$/bins (450 kilos)"
$ws.Range("R708").Value = "Región de O'Higgins"
$ws.Range("S708").Value = 289
$ws.Range("T708").Value = 450

# Make sure the D column keeps the date/datetime number format used elsewhere
# in the column (style index 2 in the original workbook).
$ws.Range("D707").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D708").NumberFormat = "YYYY-MM-DD HH:MM:SS"
